# Insert a new daily price record as row 577 of the "Zanahoria" price
# sheet. Excel shifts every existing row at/after 577 down by one
# (577->578, ..., 696->697), which matches the target diff (dimension
# grows from A1:R696 to A1:R697 and every row from 577 on now carries
# the data that used to sit one row above it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 577 (and everything below it) down by one row.
$ws.Rows.Item(577).Insert()

# Populate the newly-opened row 577 with the new record.
$ws.Range("A577").Value = 1
$ws.Range("B577").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C577").Value = "Arica y Parinacota"
$ws.Range("D577").Value = 45258
$ws.Range("E577").Value = 15
$ws.Range("F577").Value = 100114013
$ws.Range("G577").Value = "Zanahoria"
$ws.Range("H577").Value = "Sin especificar"
$ws.Range("I577").Value = "Primera"
$ws.Range("J577").Value = 60
$ws.Range("K577").Value = 13000
$ws.Range("L577").Value = 14000
$ws.Range("M577").Value = 13500
$ws.Range("N577").Value = "$/saco 25 kilos"
$ws.Range("O577").Value = "Región de Arica y Parinacota"
$ws.Range("P577").Value = 540
$ws.Range("Q577").Value = 25
$ws.Range("R577").Value = "Hortaliza"
